$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{
        r = 85
        A = "Record"
        B = "Balanço Geral"
        C = "Agricultura"
        D = "2025-04-08T13:24"
        E = "Negativo"
        F = "Moradores da localidade de Balança do Jair cobram construção de ponte. Entrevista com representante da comunidade. Moradores aguardam a entrega da ponte de concreto. A atual é de madeira e foi entregue há cerca de um ano, após enxurrada. Segundo moradora, já caiu muito gente da ponte. Em abril do ano passado, equipe esteve no local. Naquele período, situação era pior devido às fortes chuvas. moradores mesmo fizeram os reparos. Dois meses depois, medidas paliativas foram realizadas com a construção da ponte de madeira. Eles querem soluçao definitiva. Entrevista com morador e produtor rural, que falou da expectativa de ver a obra concluída. Ano passado, prefeitura prometeu esta e mais 18 pontes seriam concluídas. Entrevista com outro morador. Falta de estrutura das pontes é problema crônico. Ordem de serviço foi divulgada pela prefeitura. Processo estava em fase de assinatura do contrato, mas Tribunal de Contas não liberou a obra, por falta de licenciamento do INEA. Apresentadora disse que aguarda resposta do município, mas lembra que a resposta do município no ano passado se referia à espera pelo licenciamento. Ela cobrou do INEA o licenciamento. *matéria*"
    },
    @{
        r = 86
        A = "Record"
        B = "Balanço Geral"
        C = "Saúde"
        D = "2025-04-08T14:27"
        E = "Positivo"
        F = "Começa campanha de imunização contra a gripe para grupo prioritário. Entrevista com pessoas sendo vacinadas e com o subsecretário de Vigilância em Saúde, Charbell Kury. *matéria* também foi exibida ontem."
    },
    @{
        r = 87
        A = "Record"
        B = "Balanço Geral"
        C = "Trânsito"
        D = "2025-04-08T14:38"
        E = "Negativo"
        F = "Tô na bronca. Motociclistas reclamam da situação da Ponte Barcelos Martins. Segundo eles, o problema é na chegada ao Centro porque o semáforo está sendo coberto por galhos de árvore. Trabalho de poda já está na programação da Secretaria de Serviços Públicos. *nota coberta*"
    }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
}
